$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

$ws.Range("R4").Value = 1

$ws.Range("R22").Value = 1

$ws.Range("R23").Value = 1
$ws.Range("T23").Value = 0
$ws.Range("U23").Value = 12

$ws.Range("R25").Value = 1
$ws.Range("T25").Value = 0

$ws.Range("R38").Value = 3
$ws.Range("T38").Value = 3

$ws.Range("R39").Value = 3
$ws.Range("T39").Value = 6

$ws.Range("R46").Value = 1
$ws.Range("T46").Value = 2

$ws.Range("R47").Value = 4

$ws.Range("R52").Value = 1
$ws.Range("T52").Value = 7

$ws.Range("R54").Value = 1

$ws.Range("R63").Value = 1
$ws.Range("T63").Value = 2

$ws.Range("C71").Value = 83

# Restore auto row height on the touched rows: writing to cells on these
# rows can stamp an explicit row height, so re-run AutoFit to drop it and
# keep the row formatting identical to before the edit.
$touchedRows = @(4, 22, 23, 25, 38, 39, 46, 47, 52, 54, 63, 71)
foreach ($r in $touchedRows) {
    $ws.Rows($r).EntireRow.AutoFit()
}
